$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 91 (record #90): Warsaw LOT flight
$ws.Cells.Item(91,1).Value = 90
$ws.Cells.Item(91,2).Value = "Sunday, Jan 08"
$ws.Cells.Item(91,3).Value = "9:20 PM"
$ws.Cells.Item(91,4).Value = "LO3826"
$ws.Cells.Item(91,5).Value = "Warsaw"
$ws.Cells.Item(91,6).Value = "(WAW)"
$ws.Cells.Item(91,7).Value = "LOT "
$ws.Cells.Item(91,8).Value = "E75S"
$ws.Cells.Item(91,9).Value = "(SP-LIQ)"
$ws.Cells.Item(91,10).Value = "9:21 PM"
$ws.Cells.Item(91,11).Font.Bold = $false
$ws.Cells.Item(91,12).Value = "0 hours, 1 minutes"
$ws.Cells.Item(91,13).Font.Bold = $false

# Row 92 (record #91): Stockholm Ryanair flight
$ws.Cells.Item(92,1).Value = 91
$ws.Cells.Item(92,2).Value = "Sunday, Jan 08"
$ws.Cells.Item(92,3).Value = "9:40 PM"
$ws.Cells.Item(92,4).Value = "FR4617"
$ws.Cells.Item(92,5).Value = "Stockholm"
$ws.Cells.Item(92,6).Value = "(ARN)"
$ws.Cells.Item(92,7).Value = "Ryanair "
$ws.Cells.Item(92,8).Value = "B738"
$ws.Cells.Item(92,9).Value = "(9H-QBG)"
$ws.Cells.Item(92,10).Value = "10:26 PM"
$ws.Cells.Item(92,11).Font.Bold = $false
$ws.Cells.Item(92,12).Value = "0 hours, 46 minutes"
$ws.Cells.Item(92,13).Font.Bold = $false
